$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(57.14571807951455, 65.45588390529966, 84.98174670937232, -35.50340589055043, 57.85298202279029, 9.290402801141639, -50.55973039993673, -68.13806777921256, -22.32442281339111, -26.38550494778482, -30.35129918917852, 44.79659810853599, 53.10157927067019, 86.58216866330858, -17.10896178306069)
$row2 = @(-6.186719900796248, 8.933753480231374, 89.90615483809918, 76.47735309733278, 89.40251760147339, -37.36677148442045, 19.02214525930287, 12.50723293051084, -36.25046865710955, -4.71823884534659, 7.343408082669683, 93.07749420195671, -70.14850984229169, -18.77170567501881, -37.06698665258534)

for ($i = 0; $i -lt $row1.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $row1[$i]
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}
